$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a cell that already carries the target "N/A" style (fill/border/font)
# as the format source, so the copied cells end up on the same cellXf (style 4)
# that Excel already uses elsewhere on the sheet for "N/A" values.
$styleSource = $ws.Range("F2")
$styleSource.Copy()

$targets = @("D2:E2", "F3:G3", "D8:E8")
foreach ($addr in $targets) {
    $rng = $ws.Range($addr)
    $rng.PasteSpecial(-4122)  # xlPasteFormats
    $rng.Value = "N/A"
}

$excel.CutCopyMode = $false

# Restore the selection to D8:E8 with D8 as the active cell, matching the
# updated sheetView selection in the workbook.
$ws.Range("D8:E8").Select()
